$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 101.75
$ws.Range("I5").Value = 107.85714
$ws.Range("J5").Value = 59
$ws.Range("K5").Value = 107.85714
$ws.Range("L5").Value = 59
$ws.Range("M5").Value = 7.142859999999999
$ws.Range("N5").Value = -289

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H17").Value = 1332.3334
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1332.3334
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3997.0002
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -4333.0002

$ws.Range("H70").Value = 3205.35
$ws.Range("I70").Value = 1951.4286
$ws.Range("K70").Value = 5854.2858
$ws.Range("M70").Value = -5584.2858

$ws.Range("H73").Value = 3205.35
$ws.Range("I73").Value = 1951.4286
$ws.Range("K73").Value = 5854.2858
$ws.Range("M73").Value = -4918.2858

$ws.Range("H98").Value = 875
$ws.Range("I98").Value = 850
$ws.Range("K98").Value = 850
$ws.Range("M98").Value = 648

$ws.Range("H113").Value = 3683
$ws.Range("I113").Value = 524.5
$ws.Range("K113").Value = 524.5
$ws.Range("M113").Value = 2729.5

$ws.Range("H122").Value = 875
$ws.Range("I122").Value = 850
$ws.Range("K122").Value = 2550
$ws.Range("M122").Value = -100

$ws.Range("H132").Value = 9771.074000000001
$ws.Range("I132").Value = 8956.76
$ws.Range("J132").Value = 19950
$ws.Range("K132").Value = 26870.28
$ws.Range("L132").Value = 59850
$ws.Range("M132").Value = -24340.28
$ws.Range("N132").Value = -64910

$ws.Range("H135").Value = 920.0714
$ws.Range("I135").Value = 748.4167
$ws.Range("K135").Value = 6735.7503
$ws.Range("M135").Value = -4200.7503

$ws.Range("H137").Value = 6560.4375
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 6560.4375
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 19681.3125
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -24781.3125

$ws.Range("H138").Value = 3033
$ws.Range("J138").Value = 4904.0835
$ws.Range("L138").Value = 14712.2505
$ws.Range("N138").Value = -24992.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 90.57143000000001
$ws.Range("I5").Value = 83.5
$ws.Range("K5").Value = 83.5
$ws.Range("M5").Value = 28.5

$ws.Range("H61").Value = 7150
$ws.Range("I61").Value = 5166.6665
$ws.Range("K61").Value = 5166.6665
$ws.Range("M61").Value = -4954.6665

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("N72").ClearContents()

$ws.Range("H96").Value = 3348150
$ws.Range("J96").Value = 3348150
$ws.Range("L96").Value = 3348150
$ws.Range("N96").Value = -3353642

$ws.Range("H122").Value = 1507
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H136").Value = 7150
$ws.Range("I136").Value = 5166.6665
$ws.Range("K136").Value = 15499.9995
$ws.Range("M136").Value = -12949.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 90.57143000000001
$ws.Range("I4").Value = 83.5
$ws.Range("K4").Value = 83.5
$ws.Range("M4").Value = 31.5

$ws.Range("H20").Value = 2098.2727
$ws.Range("J20").Value = 2762.5
$ws.Range("L20").Value = 2762.5
$ws.Range("N20").Value = -3256.5

$ws.Range("H22").Value = 1450.1
$ws.Range("I22").Value = 1450.1
$ws.Range("K22").Value = 1450.1
$ws.Range("M22").Value = -1277.1

$ws.Range("H100").Value = 20020.5
$ws.Range("J100").Value = 20020.5
$ws.Range("L100").Value = 20020.5
$ws.Range("N100").Value = -22184.5

$ws.Range("H105").Value = 1838.5714
$ws.Range("I105").Value = 1728.3334
$ws.Range("K105").Value = 1728.3334
$ws.Range("M105").Value = 18.66660000000002

$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1249.5
$ws.Range("I16").Value = 1249.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1249.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -962.5
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 3749.5
$ws.Range("I22").Value = 4499
$ws.Range("K22").Value = 4499
$ws.Range("M22").Value = -4149

$ws.Range("H28").Value = 15308.4
$ws.Range("J28").Value = 15308.4
$ws.Range("L28").Value = 15308.4
$ws.Range("N28").Value = -15798.4

$ws.Range("H31").Value = 5184.3438
$ws.Range("I31").Value = 2405.353
$ws.Range("J31").Value = 8333.866
$ws.Range("K31").Value = 2405.353
$ws.Range("L31").Value = 8333.866
$ws.Range("M31").Value = -2110.353
$ws.Range("N31").Value = -8923.866

$ws.Range("H34").Value = 5184.3438
$ws.Range("I34").Value = 2405.353
$ws.Range("J34").Value = 8333.866
$ws.Range("K34").Value = 2405.353
$ws.Range("L34").Value = 8333.866
$ws.Range("M34").Value = -2203.353
$ws.Range("N34").Value = -8737.866

$ws.Range("H62").Value = 4525
$ws.Range("I62").Value = 4525
$ws.Range("K62").Value = 4525
$ws.Range("M62").Value = -3901

$ws.Range("H65").Value = 4525
$ws.Range("I65").Value = 4525
$ws.Range("K65").Value = 22625
$ws.Range("M65").Value = -19505

$ws.Range("H92").Value = 10258.429
$ws.Range("J92").Value = 10258.429
$ws.Range("L92").Value = 10258.429
$ws.Range("N92").Value = -15250.429

$ws.Range("H105").Value = 2601.1333
$ws.Range("I105").Value = 2431.1428
$ws.Range("K105").Value = 2431.1428
$ws.Range("M105").Value = -684.1428000000001

$ws.Range("H107").Value = 657.5
$ws.Range("I107").Value = 315
$ws.Range("K107").Value = 315
$ws.Range("M107").Value = 1605

$ws.Range("H113").Value = 1249.5
$ws.Range("I113").Value = 1249.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1249.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 920.5
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 2306.8572
$ws.Range("I122").Value = 1858
$ws.Range("K122").Value = 5574
$ws.Range("M122").Value = -3124

$ws.Range("H134").Value = 1732.4375
$ws.Range("I134").Value = 1732.4375
$ws.Range("K134").Value = 5197.3125
$ws.Range("M134").Value = -2662.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.2
$ws.Range("I2").Value = 26.17647
$ws.Range("J2").Value = 48.125
$ws.Range("K2").Value = 157.05882
$ws.Range("L2").Value = 288.75
$ws.Range("M2").Value = -44.05882
$ws.Range("N2").Value = -514.75

$ws.Range("H8").Value = 145.28572
$ws.Range("I8").Value = 145.28572
$ws.Range("K8").Value = 435.85716
$ws.Range("M8").Value = -296.85716

$ws.Range("H12").Value = 92.39130400000001
$ws.Range("I12").Value = 10.166667
$ws.Range("K12").Value = 30.500001
$ws.Range("M12").Value = 142.499999

$ws.Range("H18").Value = 1064.25
$ws.Range("I18").Value = 1064.25
$ws.Range("K18").Value = 3192.75
$ws.Range("M18").Value = -3023.75

$ws.Range("H52").Value = 1134.2858
$ws.Range("J52").Value = 1134.2858
$ws.Range("L52").Value = 3402.8574
$ws.Range("N52").Value = -3934.8574

$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 10000
$ws.Range("K69").Value = 30000
$ws.Range("M69").Value = -29189

$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 10000
$ws.Range("K72").Value = 90000
$ws.Range("M72").Value = -85944

$ws.Range("H81").Value = 1140.5714
$ws.Range("I81").Value = 998
$ws.Range("J81").Value = 1497
$ws.Range("K81").Value = 2994
$ws.Range("L81").Value = 4491
$ws.Range("M81").Value = -1871
$ws.Range("N81").Value = -6737

$ws.Range("H84").Value = 1140.5714
$ws.Range("I84").Value = 998
$ws.Range("J84").Value = 1497
$ws.Range("K84").Value = 8982
$ws.Range("L84").Value = 13473
$ws.Range("M84").Value = -3366
$ws.Range("N84").Value = -24705

$ws.Range("H92").Value = 692.5
$ws.Range("I92").Value = 485
$ws.Range("K92").Value = 1455
$ws.Range("M92").Value = -207

$ws.Range("H118").Value = 330
$ws.Range("I118").Value = 330
$ws.Range("K118").Value = 990
$ws.Range("M118").Value = 253

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 19999
$ws.Range("I58").Value = 19999
$ws.Range("K58").Value = 19999
$ws.Range("M58").Value = -19722

$ws.Range("H70").Value = 1500
$ws.Range("I70").Value = 1500
$ws.Range("K70").Value = 1500
$ws.Range("M70").Value = -1230

$ws.Range("H73").Value = 1500
$ws.Range("I73").Value = 1500
$ws.Range("K73").Value = 1500
$ws.Range("M73").Value = -564

$ws.Range("H80").Value = 1466.6666
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 3100
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 3100
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -5096

$ws.Range("H83").Value = 1466.6666
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 3100
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 15500
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -25484

$ws.Range("H102").Value = 1526.1177
$ws.Range("J102").Value = 2849.5
$ws.Range("L102").Value = 2849.5
$ws.Range("N102").Value = -6093.5

$ws.Range("H105").Value = 30750
$ws.Range("J105").Value = 30750
$ws.Range("L105").Value = 30750
$ws.Range("N105").Value = -37738

$ws.Range("H132").Value = 3848.75
$ws.Range("I132").Value = 1798.3334
$ws.Range("K132").Value = 5395.0002
$ws.Range("M132").Value = -2865.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4999.5
$ws.Range("I40").Value = 4999.5
$ws.Range("K40").Value = 4999.5
$ws.Range("M40").Value = -4863.5

$ws.Range("H55").Value = 1239.5
$ws.Range("I55").Value = 1239.5
$ws.Range("K55").Value = 1239.5
$ws.Range("M55").Value = -1066.5

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").ClearContents()

$ws.Range("H132").Value = 3997.5
$ws.Range("J132").Value = 3995
$ws.Range("L132").Value = 11985
$ws.Range("N132").Value = -17045

$ws.Range("H136").Value = 3106.625
$ws.Range("I136").Value = 2479.6667
$ws.Range("K136").Value = 7439.000100000001
$ws.Range("M136").Value = -4889.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 812
$ws.Range("I107").Value = 674.5
$ws.Range("J107").Value = 949.5
$ws.Range("K107").Value = 2023.5
$ws.Range("L107").Value = 2848.5
$ws.Range("M107").Value = -103.5
$ws.Range("N107").Value = -6688.5
